# Update the "想去人数" (want-to-go count) figures (column F) that changed
# between crawler runs, per the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 319
$ws.Range("F12").Value = 571
$ws.Range("F17").Value = 6630
$ws.Range("F21").Value = 7579
$ws.Range("F26").Value = 1798
$ws.Range("F29").Value = 127
$ws.Range("F32").Value = 219
$ws.Range("F34").Value = 1673
$ws.Range("F35").Value = 11
$ws.Range("F39").Value = 1204
$ws.Range("F40").Value = 1761
$ws.Range("F41").Value = 2134

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 6

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 319
$ws.Range("F12").Value = 6
$ws.Range("F14").Value = 572
$ws.Range("F20").Value = 6630
$ws.Range("F24").Value = 7579
$ws.Range("F29").Value = 1798
$ws.Range("F32").Value = 127
$ws.Range("F36").Value = 219
$ws.Range("F38").Value = 1673
$ws.Range("F39").Value = 11
$ws.Range("F44").Value = 1204
$ws.Range("F45").Value = 1761
$ws.Range("F47").Value = 2134
